# Fix Training Data Issue (#48)
#
# The NBA team box-score stats in this sheet were off by one day
# because of the way NBA stats were originally pulled/shown, so the
# numeric stat values (and their RANK columns) for every team row
# need to be corrected to the values for the intended game date, and
# the Date column (BF) needs to read "2015-05-08" instead of the
# previous mis-formatted "5-8-2014-15" label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the Date column as Text so that assigning a date-shaped
# string ("2015-05-08") is kept as literal text instead of being
# auto-converted by Excel into a date serial number.
$ws.Range("BF2:BF31").NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 45).Value = 22   # AS2: 23 -> 22
$ws.Cells.Item(2, 58).Value = "2015-05-08"   # BF2: "5-8-2014-15" -> "2015-05-08"

# Row 3
$ws.Cells.Item(3, 34).Value = 8   # AH3: 9 -> 8
$ws.Cells.Item(3, 43).Value = 14   # AQ3: 15 -> 14
$ws.Cells.Item(3, 51).Value = 21   # AY3: 22 -> 21
$ws.Cells.Item(3, 52).Value = 20   # AZ3: 21 -> 20
$ws.Cells.Item(3, 58).Value = "2015-05-08"   # BF3: "5-8-2014-15" -> "2015-05-08"

# Row 4
$ws.Cells.Item(4, 34).Value = 3   # AH4: 4 -> 3
$ws.Cells.Item(4, 46).Value = 22   # AT4: 23 -> 22
$ws.Cells.Item(4, 47).Value = 20   # AU4: 21 -> 20
$ws.Cells.Item(4, 58).Value = "2015-05-08"   # BF4: "5-8-2014-15" -> "2015-05-08"

# Row 5
$ws.Cells.Item(5, 44).Value = 26   # AR5: 25 -> 26
$ws.Cells.Item(5, 51).Value = 23   # AY5: 24 -> 23
$ws.Cells.Item(5, 58).Value = "2015-05-08"   # BF5: "5-8-2014-15" -> "2015-05-08"

# Row 6
$ws.Cells.Item(6, 32).Value = 9   # AF6: 10 -> 9
$ws.Cells.Item(6, 33).Value = 9   # AG6: 10 -> 9
$ws.Cells.Item(6, 34).Value = 8   # AH6: 9 -> 8
$ws.Cells.Item(6, 48).Value = 12   # AV6: 13 -> 12
$ws.Cells.Item(6, 58).Value = "2015-05-08"   # BF6: "5-8-2014-15" -> "2015-05-08"

# Row 7
$ws.Cells.Item(7, 34).Value = 23   # AH7: 24 -> 23
$ws.Cells.Item(7, 42).Value = 12   # AP7: 11 -> 12
$ws.Cells.Item(7, 43).Value = 18   # AQ7: 17 -> 18
$ws.Cells.Item(7, 45).Value = 21   # AS7: 22 -> 21
$ws.Cells.Item(7, 58).Value = "2015-05-08"   # BF7: "5-8-2014-15" -> "2015-05-08"

# Row 8
$ws.Cells.Item(8, 4).Value = 82   # D8: 81 -> 82
$ws.Cells.Item(8, 6).Value = 32   # F8: 31 -> 32
$ws.Cells.Item(8, 7).Value = 0.61   # G8: 0.617 -> 0.61
$ws.Cells.Item(8, 10).Value = 85.8   # J8: 85.90000000000001 -> 85.8
$ws.Cells.Item(8, 11).Value = 0.463   # K8: 0.462 -> 0.463
$ws.Cells.Item(8, 14).Value = 0.352   # N8: 0.351 -> 0.352
$ws.Cells.Item(8, 17).Value = 0.752   # Q8: 0.751 -> 0.752
$ws.Cells.Item(8, 19).Value = 31.8   # S8: 31.9 -> 31.8
$ws.Cells.Item(8, 20).Value = 42.3   # T8: 42.4 -> 42.3
$ws.Cells.Item(8, 21).Value = 22.5   # U8: 22.6 -> 22.5
$ws.Cells.Item(8, 27).Value = 22.1   # AA8: 22.2 -> 22.1
$ws.Cells.Item(8, 28).Value = 105.2   # AB8: 105.3 -> 105.2
$ws.Cells.Item(8, 29).Value = 2.9   # AC8: 3 -> 2.9
$ws.Cells.Item(8, 30).Value = 1   # AD8: 25 -> 1
$ws.Cells.Item(8, 32).Value = 9   # AF8: 8 -> 9
$ws.Cells.Item(8, 36).Value = 8   # AJ8: 6 -> 8
$ws.Cells.Item(8, 43).Value = 16   # AQ8: 18 -> 16
$ws.Cells.Item(8, 45).Value = 23   # AS8: 21 -> 23
$ws.Cells.Item(8, 46).Value = 23   # AT8: 22 -> 23
$ws.Cells.Item(8, 50).Value = 20   # AX8: 18 -> 20
$ws.Cells.Item(8, 58).Value = "2015-05-08"   # BF8: "5-8-2014-15" -> "2015-05-08"

# Row 9
$ws.Cells.Item(9, 34).Value = 8   # AH9: 9 -> 8
$ws.Cells.Item(9, 45).Value = 12   # AS9: 13 -> 12
$ws.Cells.Item(9, 50).Value = 21   # AX9: 22 -> 21
$ws.Cells.Item(9, 58).Value = "2015-05-08"   # BF9: "5-8-2014-15" -> "2015-05-08"

# Row 10
$ws.Cells.Item(10, 36).Value = 6   # AJ10: 7 -> 6
$ws.Cells.Item(10, 53).Value = 24   # BA10: 23 -> 24
$ws.Cells.Item(10, 54).Value = 18   # BB10: 19 -> 18
$ws.Cells.Item(10, 58).Value = "2015-05-08"   # BF10: "5-8-2014-15" -> "2015-05-08"

# Row 11
$ws.Cells.Item(11, 58).Value = "2015-05-08"   # BF11: "5-8-2014-15" -> "2015-05-08"

# Row 12
$ws.Cells.Item(12, 4).Value = 82   # D12: 81 -> 82
$ws.Cells.Item(12, 5).Value = 56   # E12: 55 -> 56
$ws.Cells.Item(12, 7).Value = 0.6830000000000001   # G12: 0.679 -> 0.6830000000000001
$ws.Cells.Item(12, 10).Value = 83.3   # J12: 83.40000000000001 -> 83.3
$ws.Cells.Item(12, 14).Value = 0.348   # N12: 0.347 -> 0.348
$ws.Cells.Item(12, 15).Value = 18.6   # O12: 18.4 -> 18.6
$ws.Cells.Item(12, 16).Value = 26   # P12: 25.7 -> 26
$ws.Cells.Item(12, 17).Value = 0.715   # Q12: 0.716 -> 0.715
$ws.Cells.Item(12, 20).Value = 43.7   # T12: 43.6 -> 43.7
$ws.Cells.Item(12, 25).Value = 5.3   # Y12: 5.4 -> 5.3
$ws.Cells.Item(12, 26).Value = 22   # Z12: 21.9 -> 22
$ws.Cells.Item(12, 27).Value = 21.1   # AA12: 21 -> 21.1
$ws.Cells.Item(12, 29).Value = 3.4   # AC12: 3.3 -> 3.4
$ws.Cells.Item(12, 30).Value = 1   # AD12: 25 -> 1
$ws.Cells.Item(12, 31).Value = 3   # AE12: 4 -> 3
$ws.Cells.Item(12, 33).Value = 3   # AG12: 4 -> 3
$ws.Cells.Item(12, 41).Value = 5   # AO12: 6 -> 5
$ws.Cells.Item(12, 42).Value = 2   # AP12: 3 -> 2
$ws.Cells.Item(12, 46).Value = 14   # AT12: 15 -> 14
$ws.Cells.Item(12, 51).Value = 22   # AY12: 23 -> 22
$ws.Cells.Item(12, 58).Value = "2015-05-08"   # BF12: "5-8-2014-15" -> "2015-05-08"

# Row 13
$ws.Cells.Item(13, 43).Value = 13   # AQ13: 14 -> 13
$ws.Cells.Item(13, 58).Value = "2015-05-08"   # BF13: "5-8-2014-15" -> "2015-05-08"

# Row 14
$ws.Cells.Item(14, 58).Value = "2015-05-08"   # BF14: "5-8-2014-15" -> "2015-05-08"

# Row 15
$ws.Cells.Item(15, 4).Value = 82   # D15: 81 -> 82
$ws.Cells.Item(15, 6).Value = 61   # F15: 60 -> 61
$ws.Cells.Item(15, 7).Value = 0.256   # G15: 0.259 -> 0.256
$ws.Cells.Item(15, 9).Value = 37.2   # I15: 37.4 -> 37.2
$ws.Cells.Item(15, 10).Value = 85.59999999999999   # J15: 85.7 -> 85.59999999999999
$ws.Cells.Item(15, 11).Value = 0.435   # K15: 0.436 -> 0.435
$ws.Cells.Item(15, 13).Value = 18.9   # M15: 19 -> 18.9
$ws.Cells.Item(15, 15).Value = 17.5   # O15: 17.3 -> 17.5
$ws.Cells.Item(15, 16).Value = 23.6   # P15: 23.4 -> 23.6
$ws.Cells.Item(15, 17).Value = 0.741   # Q15: 0.739 -> 0.741
$ws.Cells.Item(15, 19).Value = 32.3   # S15: 32.4 -> 32.3
$ws.Cells.Item(15, 20).Value = 43.9   # T15: 44 -> 43.9
$ws.Cells.Item(15, 21).Value = 20.9   # U15: 21 -> 20.9
$ws.Cells.Item(15, 25).Value = 4.8   # Y15: 4.9 -> 4.8
$ws.Cells.Item(15, 26).Value = 21.2   # Z15: 21.1 -> 21.2
$ws.Cells.Item(15, 27).Value = 19.4   # AA15: 19.2 -> 19.4
$ws.Cells.Item(15, 28).Value = 98.5   # AB15: 98.59999999999999 -> 98.5
$ws.Cells.Item(15, 29).Value = -6.8   # AC15: -6.7 -> -6.8
$ws.Cells.Item(15, 30).Value = 1   # AD15: 25 -> 1
$ws.Cells.Item(15, 42).Value = 11   # AP15: 13 -> 11
$ws.Cells.Item(15, 45).Value = 13   # AS15: 12 -> 13
$ws.Cells.Item(15, 46).Value = 12   # AT15: 11 -> 12
$ws.Cells.Item(15, 47).Value = 21   # AU15: 20 -> 21
$ws.Cells.Item(15, 50).Value = 22   # AX15: 21 -> 22
$ws.Cells.Item(15, 52).Value = 21   # AZ15: 20 -> 21
$ws.Cells.Item(15, 53).Value = 23   # BA15: 24 -> 23
$ws.Cells.Item(15, 54).Value = 19   # BB15: 17 -> 19
$ws.Cells.Item(15, 58).Value = "2015-05-08"   # BF15: "5-8-2014-15" -> "2015-05-08"

# Row 16
$ws.Cells.Item(16, 31).Value = 5   # AE16: 4 -> 5
$ws.Cells.Item(16, 37).Value = 9   # AK16: 10 -> 9
$ws.Cells.Item(16, 58).Value = "2015-05-08"   # BF16: "5-8-2014-15" -> "2015-05-08"

# Row 17
$ws.Cells.Item(17, 38).Value = 21   # AL17: 22 -> 21
$ws.Cells.Item(17, 50).Value = 18   # AX17: 19 -> 18
$ws.Cells.Item(17, 58).Value = "2015-05-08"   # BF17: "5-8-2014-15" -> "2015-05-08"

# Row 18
$ws.Cells.Item(18, 34).Value = 3   # AH18: 4 -> 3
$ws.Cells.Item(18, 58).Value = "2015-05-08"   # BF18: "5-8-2014-15" -> "2015-05-08"

# Row 19
$ws.Cells.Item(19, 34).Value = 20   # AH19: 21 -> 20
$ws.Cells.Item(19, 42).Value = 3   # AP19: 2 -> 3
$ws.Cells.Item(19, 58).Value = "2015-05-08"   # BF19: "5-8-2014-15" -> "2015-05-08"

# Row 20
$ws.Cells.Item(20, 4).Value = 82   # D20: 81 -> 82
$ws.Cells.Item(20, 5).Value = 45   # E20: 44 -> 45
$ws.Cells.Item(20, 7).Value = 0.549   # G20: 0.543 -> 0.549
$ws.Cells.Item(20, 10).Value = 82.90000000000001   # J20: 82.59999999999999 -> 82.90000000000001
$ws.Cells.Item(20, 11).Value = 0.457   # K20: 0.458 -> 0.457
$ws.Cells.Item(20, 12).Value = 7.1   # L20: 7.2 -> 7.1
$ws.Cells.Item(20, 14).Value = 0.37   # N20: 0.372 -> 0.37
$ws.Cells.Item(20, 16).Value = 21.8   # P20: 21.7 -> 21.8
$ws.Cells.Item(20, 17).Value = 0.751   # Q20: 0.756 -> 0.751
$ws.Cells.Item(20, 18).Value = 11.5   # R20: 11.3 -> 11.5
$ws.Cells.Item(20, 19).Value = 32   # S20: 31.9 -> 32
$ws.Cells.Item(20, 20).Value = 43.5   # T20: 43.2 -> 43.5
$ws.Cells.Item(20, 24).Value = 6.2   # X20: 6.1 -> 6.2
$ws.Cells.Item(20, 27).Value = 18.7   # AA20: 18.6 -> 18.7
$ws.Cells.Item(20, 28).Value = 99.40000000000001   # AB20: 99.3 -> 99.40000000000001
$ws.Cells.Item(20, 29).Value = 0.8   # AC20: 0.6 -> 0.8
$ws.Cells.Item(20, 30).Value = 1   # AD20: 25 -> 1
$ws.Cells.Item(20, 31).Value = 13   # AE20: 14 -> 13
$ws.Cells.Item(20, 33).Value = 13   # AG20: 14 -> 13
$ws.Cells.Item(20, 35).Value = 11   # AI20: 12 -> 11
$ws.Cells.Item(20, 36).Value = 19   # AJ20: 22 -> 19
$ws.Cells.Item(20, 37).Value = 10   # AK20: 9 -> 10
$ws.Cells.Item(20, 43).Value = 17   # AQ20: 13 -> 17
$ws.Cells.Item(20, 45).Value = 19   # AS20: 20 -> 19
$ws.Cells.Item(20, 49).Value = 25   # AW20: 26 -> 25
$ws.Cells.Item(20, 55).Value = 13   # BC20: 14 -> 13
$ws.Cells.Item(20, 58).Value = "2015-05-08"   # BF20: "5-8-2014-15" -> "2015-05-08"

# Row 21
$ws.Cells.Item(21, 58).Value = "2015-05-08"   # BF21: "5-8-2014-15" -> "2015-05-08"

# Row 22
$ws.Cells.Item(22, 41).Value = 6   # AO22: 5 -> 6
$ws.Cells.Item(22, 43).Value = 15   # AQ22: 16 -> 15
$ws.Cells.Item(22, 58).Value = "2015-05-08"   # BF22: "5-8-2014-15" -> "2015-05-08"

# Row 23
$ws.Cells.Item(23, 4).Value = 82   # D23: 81 -> 82
$ws.Cells.Item(23, 6).Value = 57   # F23: 56 -> 57
$ws.Cells.Item(23, 7).Value = 0.305   # G23: 0.309 -> 0.305
$ws.Cells.Item(23, 9).Value = 37.5   # I23: 37.6 -> 37.5
$ws.Cells.Item(23, 11).Value = 0.453   # K23: 0.454 -> 0.453
$ws.Cells.Item(23, 13).Value = 19.5   # M23: 19.6 -> 19.5
$ws.Cells.Item(23, 19).Value = 31.8   # S23: 31.7 -> 31.8
$ws.Cells.Item(23, 20).Value = 41.8   # T23: 41.6 -> 41.8
$ws.Cells.Item(23, 21).Value = 20.6   # U23: 20.7 -> 20.6
$ws.Cells.Item(23, 25).Value = 5.4   # Y23: 5.2 -> 5.4
$ws.Cells.Item(23, 28).Value = 95.7   # AB23: 95.8 -> 95.7
$ws.Cells.Item(23, 29).Value = -5.7   # AC23: -5.6 -> -5.7
$ws.Cells.Item(23, 30).Value = 1   # AD23: 25 -> 1
$ws.Cells.Item(23, 36).Value = 20   # AJ23: 19 -> 20
$ws.Cells.Item(23, 38).Value = 22   # AL23: 21 -> 22
$ws.Cells.Item(23, 44).Value = 25   # AR23: 26 -> 25
$ws.Cells.Item(23, 51).Value = 24   # AY23: 21 -> 24
$ws.Cells.Item(23, 52).Value = 19   # AZ23: 18 -> 19
$ws.Cells.Item(23, 58).Value = "2015-05-08"   # BF23: "5-8-2014-15" -> "2015-05-08"

# Row 24
$ws.Cells.Item(24, 34).Value = 20   # AH24: 21 -> 20
$ws.Cells.Item(24, 36).Value = 22   # AJ24: 21 -> 22
$ws.Cells.Item(24, 58).Value = "2015-05-08"   # BF24: "5-8-2014-15" -> "2015-05-08"

# Row 25
$ws.Cells.Item(25, 34).Value = 8   # AH25: 9 -> 8
$ws.Cells.Item(25, 36).Value = 6   # AJ25: 7 -> 6
$ws.Cells.Item(25, 58).Value = "2015-05-08"   # BF25: "5-8-2014-15" -> "2015-05-08"

# Row 26
$ws.Cells.Item(26, 50).Value = 18   # AX26: 19 -> 18
$ws.Cells.Item(26, 58).Value = "2015-05-08"   # BF26: "5-8-2014-15" -> "2015-05-08"

# Row 27
$ws.Cells.Item(27, 49).Value = 26   # AW27: 25 -> 26
$ws.Cells.Item(27, 58).Value = "2015-05-08"   # BF27: "5-8-2014-15" -> "2015-05-08"

# Row 28
$ws.Cells.Item(28, 4).Value = 82   # D28: 81 -> 82
$ws.Cells.Item(28, 5).Value = 55   # E28: 54 -> 55
$ws.Cells.Item(28, 7).Value = 0.671   # G28: 0.667 -> 0.671
$ws.Cells.Item(28, 10).Value = 83.59999999999999   # J28: 83.8 -> 83.59999999999999
$ws.Cells.Item(28, 11).Value = 0.468   # K28: 0.467 -> 0.468
$ws.Cells.Item(28, 12).Value = 8.300000000000001   # L28: 8.199999999999999 -> 8.300000000000001
$ws.Cells.Item(28, 14).Value = 0.367   # N28: 0.364 -> 0.367
$ws.Cells.Item(28, 16).Value = 21.4   # P28: 21.5 -> 21.4
$ws.Cells.Item(28, 19).Value = 33.8   # S28: 33.9 -> 33.8
$ws.Cells.Item(28, 20).Value = 43.6   # T28: 43.7 -> 43.6
$ws.Cells.Item(28, 22).Value = 14   # V28: 13.9 -> 14
$ws.Cells.Item(28, 29).Value = 6.2   # AC28: 6.3 -> 6.2
$ws.Cells.Item(28, 30).Value = 1   # AD28: 25 -> 1
$ws.Cells.Item(28, 31).Value = 5   # AE28: 6 -> 5
$ws.Cells.Item(28, 33).Value = 5   # AG28: 6 -> 5
$ws.Cells.Item(28, 42).Value = 23   # AP28: 22 -> 23
$ws.Cells.Item(28, 46).Value = 15   # AT28: 14 -> 15
$ws.Cells.Item(28, 48).Value = 13   # AV28: 12 -> 13
$ws.Cells.Item(28, 58).Value = "2015-05-08"   # BF28: "5-8-2014-15" -> "2015-05-08"

# Row 29
$ws.Cells.Item(29, 52).Value = 18   # AZ29: 19 -> 18
$ws.Cells.Item(29, 58).Value = "2015-05-08"   # BF29: "5-8-2014-15" -> "2015-05-08"

# Row 30
$ws.Cells.Item(30, 42).Value = 13   # AP30: 12 -> 13
$ws.Cells.Item(30, 45).Value = 20   # AS30: 19 -> 20
$ws.Cells.Item(30, 46).Value = 11   # AT30: 12 -> 11
$ws.Cells.Item(30, 58).Value = "2015-05-08"   # BF30: "5-8-2014-15" -> "2015-05-08"

# Row 31
$ws.Cells.Item(31, 34).Value = 1   # AH31: 2 -> 1
$ws.Cells.Item(31, 36).Value = 21   # AJ31: 20 -> 21
$ws.Cells.Item(31, 42).Value = 22   # AP31: 23 -> 22
$ws.Cells.Item(31, 54).Value = 17   # BB31: 18 -> 17
$ws.Cells.Item(31, 55).Value = 14   # BC31: 13 -> 14
$ws.Cells.Item(31, 58).Value = "2015-05-08"   # BF31: "5-8-2014-15" -> "2015-05-08"
